$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking price cells so Excel keeps them as text
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "73.440.71"
$ws.Range("E2").Value = "  +1.98%  "

# Row 3
$ws.Range("D3").Value = "4.062.22"
$ws.Range("E3").Value = "  +1.34%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").Value = "572.56"
$ws.Range("E5").Value = "  +7.17%  "

# Row 6
$ws.Range("D6").Value = "151.74"
$ws.Range("E6").Value = "  -0.52%  "

# Row 7
$ws.Range("D7").Value = "4.055.92"
$ws.Range("E7").Value = "  +1.39%  "

# Row 8
$ws.Range("D8").Value = "0.696"
$ws.Range("E8").Value = "  +0.02%  "

# Row 9
$ws.Range("E9").Value = "  -0.05%  "

# Row 10
$ws.Range("D10").Value = "0.766"
$ws.Range("E10").Value = "  +2.58%  "

# Row 11
$ws.Range("E11").Value = "  +0.67%  "

# Row 12
$ws.Range("D12").Value = "54.11"

# Row 13
$ws.Range("D13").Value = "0.0000329"
$ws.Range("E13").Value = "  +1.12%  "

# Row 14
$ws.Range("D14").Value = "11.24"
$ws.Range("E14").Value = "  +5.63%  "

# Row 15
$ws.Range("D15").Value = "4.707.31"
$ws.Range("E15").Value = "  +1.23%  "

# Row 16
$ws.Range("D16").Value = "4.056.93"
$ws.Range("E16").Value = "  +1.33%  "

# Row 17
$ws.Range("D17").Value = "14.41"
$ws.Range("E17").Value = "  +3.45%  "

# Row 18
$ws.Range("E18").Value = "  +2.04%  "

# Row 19
$ws.Range("E19").Value = "  +3.28%  "

# Row 20
$ws.Range("E20").Value = "  -0.07%  "

# Row 21
$ws.Range("D21").Value = "73.336.03"
$ws.Range("E21").Value = "  +1.97%  "

# Row 22
$ws.Range("D22").Value = "446.11"
$ws.Range("E22").Value = "  +4.69%  "

# Row 23
$ws.Range("D23").Value = "4.57"
$ws.Range("E23").Value = "  +8.57%  "

# Row 24
$ws.Range("D24").Value = "98.66"
$ws.Range("E24").Value = "  +0.86%  "

# Row 25
$ws.Range("D25").Value = "3.60"
$ws.Range("E25").Value = "  +2.94%  "

# Row 26
$ws.Range("E26").Value = "  +2.93%  "

# Row 27
$ws.Range("D27").Value = "4.28"
$ws.Range("E27").Value = "  +18.33%  "

# Row 28
$ws.Range("D28").Value = "11.51"
$ws.Range("E28").Value = "  +3.19%  "

# Row 29
$ws.Range("E29").Value = "  +4.36%  "

# Row 30
$ws.Range("E30").Value = "  +2.03%  "

# Row 31
$ws.Range("D31").Value = "37.34"
$ws.Range("E31").Value = "  +1.76%  "

# Row 32
$ws.Range("D32").Value = "7.86"
$ws.Range("E32").Value = "  +10.43%  "

# Row 33
$ws.Range("E33").Value = "  +4.07%  "

# Row 34
$ws.Range("D34").Value = "13.68"
$ws.Range("E34").Value = "  +2.43%  "

# Row 35
$ws.Range("D35").Value = "687.13"
$ws.Range("E35").Value = "  +1.97%  "

# Row 36
$ws.Range("D36").Value = "48.54"
$ws.Range("E36").Value = "  +13.90%  "

# Row 37
$ws.Range("D37").Value = "68.22"
$ws.Range("E37").Value = "  +3.75%  "

# Row 38
$ws.Range("D38").Value = "0.0₃0908"
$ws.Range("E38").Value = "  +9.90%  "

# Row 39
$ws.Range("E39").Value = "  +5.04%  "

# Row 40
$ws.Range("D40").Value = "0.150"
$ws.Range("E40").Value = "  -1.40%  "

# Row 41
$ws.Range("B41").Value = "THORChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D41").Value = "11.31"
$ws.Range("E41").Value = "  +16.73%  "

# Row 42
$ws.Range("B42").Value = "ThetaToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D42").Value = "3.39"
$ws.Range("E42").Value = "  -1.69%  "

# Row 43
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  -0.01%  "

# Row 44
$ws.Range("E44").Value = "  +1.40%  "

# Row 45
$ws.Range("E45").Value = "  +1.97%  "

# Row 46
$ws.Range("D46").Value = "0.999"
$ws.Range("E46").Value = "  +0.08%  "

# Row 47
$ws.Range("E47").Value = "  +1.57%  "

# Row 48
$ws.Range("E48").Value = "  +5.56%  "

# Row 49
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").Value = "2.24"
$ws.Range("E49").Value = "  +12.05%  "

# Row 50
$ws.Range("B50").Value = "LidoDAOToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D50").Value = "3.52"
$ws.Range("E50").Value = "  +7.43%  "

# Row 51
$ws.Range("B51").Value = "ApeXProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D51").Value = "3.33"
$ws.Range("E51").Value = "  -0.95%  "
